$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.502.82"
$ws.Range("D3").Value = "2.339.61"
$ws.Range("E3").Value = "  -1.51%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'304.86"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("D6").Value = "'101.06"
$ws.Range("E6").Value = "  -3.57%  "
$ws.Range("D7").Value = "'0.508"
$ws.Range("E7").Value = "  -3.82%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.507"
$ws.Range("E9").Value = "  -3.03%  "
$ws.Range("D10").Value = "'35.13"
$ws.Range("E10").Value = "  -3.39%  "
$ws.Range("D11").Value = "'0.0798"
$ws.Range("E11").Value = "  -2.05%  "
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "'6.80"
$ws.Range("E13").Value = "  -2.97%  "
$ws.Range("D14").Value = "2.706.77"
$ws.Range("E14").Value = "  -1.32%  "
$ws.Range("D15").Value = "'15.68"
$ws.Range("E15").Value = "  -0.37%  "
$ws.Range("D16").Value = "2.308.89"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("E17").Value = "  -0.98%  "
$ws.Range("D18").Value = "43.408.77"
$ws.Range("E19").Value = "  -1.71%  "
$ws.Range("D20").Value = "0.0₃0909"
$ws.Range("E20").Value = "  -2.22%  "
$ws.Range("D21").Value = "'6.11"
$ws.Range("E21").Value = "  -3.25%  "
$ws.Range("D22").Value = "'68.18"
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").Value = "'237.78"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "'1.98"
$ws.Range("E24").Value = "  -4.45%  "
$ws.Range("D25").Value = "'2.53"
$ws.Range("E25").Value = "  -3.71%  "
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("D27").Value = "'25.08"
$ws.Range("E27").Value = "  -3.45%  "
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'34.59"
$ws.Range("E28").Value = "  -5.54%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").Value = "'2.09"
$ws.Range("E29").Value = "  -5.49%  "
$ws.Range("D30").Value = "'165.99"
$ws.Range("E30").Value = "  +2.28%  "
$ws.Range("D31").Value = "'9.24"
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'5.05"
$ws.Range("E33").Value = "  -4.80%  "
$ws.Range("E34").Value = "  -4.77%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("D36").Value = "'16.86"
$ws.Range("E36").Value = "  -8.11%  "
$ws.Range("D37").Value = "'0.0705"
$ws.Range("E37").Value = "  -5.07%  "
$ws.Range("D38").Value = "'2.91"
$ws.Range("E38").Value = "  -7.68%  "
$ws.Range("D39").Value = "'1.82"
$ws.Range("E39").Value = "  -6.35%  "
$ws.Range("D40").Value = "'0.102"
$ws.Range("E40").Value = "  -3.69%  "
$ws.Range("E41").Value = "  -3.39%  "
$ws.Range("D42").Value = "'2.42"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "1.981.95"
$ws.Range("E43").Value = "  -1.48%  "
$ws.Range("E44").Value = "  -2.38%  "
$ws.Range("D45").Value = "'18.47"
$ws.Range("E45").Value = "  -9.86%  "
$ws.Range("D46").Value = "'10.01"
$ws.Range("E46").Value = "  -4.14%  "
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "  -7.42%  "
$ws.Range("D48").Value = "'56.37"
$ws.Range("E48").Value = "  -3.18%  "
$ws.Range("D49").Value = "'4.85"
$ws.Range("E49").Value = "  +2.43%  "
$ws.Range("D50").Value = "2.566.28"
$ws.Range("E50").Value = "  +0.21%  "
$ws.Range("D51").Value = "'1.54"
$ws.Range("E51").Value = "  -1.95%  "
